# The "derive all to railway" scenario (columns D:E) is being dropped entirely,
# leaving only the "current situation" scenario (columns B:C). Column B becomes
# "Railway" and column C becomes "Roadway" in the header row 2, and a handful of
# numeric data points in column B are corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("global_results")

# Delete columns D and E entirely (the "derive all to railway" scenario block).
# This also removes the now-unused "derive all to railway" shared string.
$ws.Range("D1:E17").Delete(-4159) | Out-Null

# Row 2 headers for the remaining "current situation" block.
$ws.Range("B2").Value = "Railway"
$ws.Range("C2").Value = "Roadway"

# Corrected numeric values (column B / Railway).
$ws.Range("B3").Value = 0.01812407162286442
$ws.Range("B4").Value = 0.0276558058426621
$ws.Range("B5").Value = 0.008085601613934908
$ws.Range("B6").Value = 0.05386547907946142
$ws.Range("B8").Value = 26412937.54733584
$ws.Range("B9").Value = 11260463543.3734
$ws.Range("B10").Value = 606550263.4206179
$ws.Range("B14").Value = 30.55844680252241
$ws.Range("B15").Value = 426.3237863336103
$ws.Range("B16").Value = 17278.3
$ws.Range("B17").Value = 651711.310914465
